$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = 5.2011951
$ws.Range("O2").Value = 100.4944918

$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = 5.2011951
$ws.Range("O4").Value = 100.4944918

$ws.Range("N28").Value = 1.4888923
$ws.Range("O28").Value = 103.7957151
